$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A3").Value = 'SC-5,SC-5 (2),CM-6 b'
$ws.Range("A4").Value = 'AC-6 (8),AC-6 (9),AU-12 (3),AU-7 a,AU-7 b,CM-5 (1),AU-8 b'
$ws.Range("A5").Value = 'CM-7 b,CM-6 b,AC-17 (9),AC-17 (1)'
$ws.Range("A15").Value = 'AU-3 (1),IA-8,IA-2'
$ws.Range("A17").Value = 'AU-12 c,AU-12 a,AU-3,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A19").Value = 'IA-5 (1) (a),IA-5 (1) (b),CM-6 b'
$ws.Range("A21").Value = 'AC-12,MA-4 (7),MA-4 e,SC-10'
$ws.Range("A22").Value = 'AU-7 (1),AU-14 (1),AU-12 a,AU-6 (4),AU-3,AU-7 a,AU-3 (1),CM-5 (1),MA-4 (1) (a),CM-6 b'
$ws.Range("A25").Value = 'AU-12 c,AU-12 a,AU-3,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A29").Value = 'SC-8 (1),SC-8,SC-8 (2)'
$ws.Range("A31").Value = 'AU-12 c,AU-12 a,AU-3,AU-3 (1),MA-4 (1) (a),AC-2 (4)'
$ws.Range("A38").Value = 'AU-9,SI-11 b'
$ws.Range("A39").Value = 'CM-6 b,AU-3'
$ws.Range("A45").Value = 'AC-8 c 1, AC-8 c 2, AC-8 c 3,AC-8 a,AC-8 b'
$ws.Range("A53").Value = 'MA-4 (6),SC-13'
$ws.Range("A55").Value = 'SC-8,AC-17 (2)'
$ws.Range("A63").Value = 'AU-5 (1),AU-5 a'
$ws.Range("A65").Value = 'CM-6 b,IA-2 (2)'
$ws.Range("A67").Value = 'AU-12 c,AU-12 a,AU-3,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A69").Value = 'AU-12 c,AU-12 a,AU-7 a,AU-12 (3),AU-7 b,CM-5 (1),AU-8 b,CM-6 b'
$ws.Range("A71").Value = 'AU-4 (1),AU-3'
$ws.Range("A77").Value = 'AU-12 c,AU-12 a,AU-3,AU-3 (1),MA-4 (1) (a),AC-2 (4)'
$ws.Range("A79").Value = 'AU-9,AU-9 (3)'
$ws.Range("A80").Value = 'IA-2 (3),IA-2 (1),IA-2 (4),IA-2 (2)'
$ws.Range("A86").Value = 'AU-12 c,AU-12 a,AU-3,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A88").Value = 'AC-2 (4),AC-6 (9),AU-12 c,CM-5 (1)'
$ws.Range("A89").Value = 'IA-2 (3),IA-2 (5),IA-2 (4),IA-2,IA-2 (2)'
$ws.Range("A91").Value = 'AU-9,AU-9 (3)'
$ws.Range("A96").Value = 'SC-8 (1),SC-8,AC-18 (1)'
$ws.Range("A97").Value = 'AU-8 b,AU-8 (1) (a),AU-8 (1) (b)'
$ws.Range("A101").Value = 'IA-11,AC-3 (4)'
$ws.Range("A102").Value = 'AU-12 c,AU-12 a,AU-3,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A111").Value = 'AU-5 b,AU-5 a'
$ws.Range("A119").Value = 'AU-12 c,AU-12 a,AU-3,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A123").Value = 'CM-7 b,CM-7 a'
$ws.Range("A124").Value = 'AU-12 c,AU-12 a,AU-3,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A128").Value = 'CM-6 b,CM-7 a,IA-5 (1) (c)'
$ws.Range("A139").Value = 'SI-6 d,CM-3 (5),SI-6 b'
$ws.Range("A148").Value = 'AU-12 c,AU-14 (1),AU-12 a,AU-3,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A157").Value = 'AU-12 c,AU-12 a,AU-3,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A159").Value = 'SC-8,AC-17 (2)'
$ws.Range("A175").Value = 'SI-16,CM-7 a'
